$d = $word.ActiveDocument

# 1. Replace the whole "green tiles" paragraph text with the updated copy.
#    This also removes the proofErr gramStart/gramEnd markers around "500"
#    since Find/Replace only rewrites run text, not structural proofErr
#    elements (those get dropped because the runs they flank get merged
#    into a fresh replacement run).
$oldText = "In the game, the green tiles on the display represent land. On the green tiles, you can have cities, bakeries, space machines, and even witchcraft. The initial costs are 500 for a city, 250 for a bakery, 4,00 for a space machine, and 100,000 for witchcraft. The costs will go up as the player purchases more buildings. The units you can have on tiles are soccer moms, soccer mom cavalries, famers, and rockets. The initial costs are 50 for a soccer mom, 150 for a soccer mom cavalry, 500 for a farmer, and 40,000 for a rocket. The costs will increase after the player purchases their initial unit. The soccer moms have a high offense, but very low mobility. On the other hand, the soccer mom cavalries have a medium offense, but have high mobility. The farmers have a high defense and good range. Finally, the rockets have a high range and they are also endgame units. "

$newText = "In the game, the green tiles on the display represent land. On the green tiles, you can have cities, bakeries, space machines, and even witchcraft. The initial costs are 500 for a city, 250 for a bakery, 4,000 for a space machine, and 15,000 for witchcraft. The costs will go up as the player purchases more buildings. The units you can have on tiles are soccer moms, soccer mom cavalries, famers, and rockets. The initial costs are 50 for a soccer mom, 250 for a soccer mom cavalry, 750 for a farmer, and 5,000 for a rocket. The costs will increase after the player purchases their initial unit. The soccer moms have a high offense, but very low mobility. On the other hand, the soccer mom cavalries have a medium offense, but have high mobility. The farmers have a high defense. Finally, the rockets have a high range and offense. "

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null

# 2. Move the "_GoBack" bookmark from the end of the document (after
#    "Graphic Designer - Erik Law") to inside the paragraph we just edited,
#    right after "In the game, the green tile".
foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

$findRange = $d.Content.Duplicate
$findRange.Find.Execute("In the game, the green tile", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmStart = $findRange.End
$d.Bookmarks.Add("_GoBack", $d.Range($bmStart, $bmStart)) | Out-Null

# 3. Update section page setup: add header/footer distances and an
#    explicit column spacing value. Word's PageSetup distances are in
#    points (1 pt = 20 twips), and the OOXML attributes (w:header /
#    w:footer) are in twips, so 720 twips = 36 pt.
$sec = $d.Sections.Item(1)
$sec.PageSetup.HeaderDistance = 36
$sec.PageSetup.FooterDistance = 36
$sec.PageSetup.TextColumns.Spacing = 36
